$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 285, shifting existing rows 285:394 down to 286:395
$ws.Rows(285).Insert()

# Populate the new row 285 with the new weekly record.
# Columns that stay constant across every record in this sheet (A,B,C,E,F,G,H,N,Q,R)
$ws.Cells.Item(285, 1).Value = 5
$ws.Cells.Item(285, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(285, 3).Value = "Maule"
$ws.Cells.Item(285, 4).Value = 45009
$ws.Cells.Item(285, 5).Value = 7
$ws.Cells.Item(285, 6).Value = 100112008
$ws.Cells.Item(285, 7).Value = "Coliflor"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 5000
$ws.Cells.Item(285, 11).Value = 800
$ws.Cells.Item(285, 12).Value = 900
$ws.Cells.Item(285, 13).Value = 860
$ws.Cells.Item(285, 14).Value = "$/unidad"
$ws.Cells.Item(285, 15).Value = "Región del Maule"
$ws.Cells.Item(285, 16).Value = 860
$ws.Cells.Item(285, 17).Value = 1
$ws.Cells.Item(285, 18).Value = "Hortaliza"
